# Update cryptocurrency price/volume figures per the latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.618.24"
$ws.Range("E2").Value = "  -1.83%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.587.91"
$ws.Range("E3").Value = "  -2.28%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.93"
$ws.Range("E5").Value = "  -1.44%  "
$ws.Range("E6").Value = "  -2.49%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  -2.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.58"
$ws.Range("E10").Value = "  -3.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0834"
$ws.Range("E11").Value = "  -1.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.810.25"
$ws.Range("E12").Value = "  -2.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.587.25"
$ws.Range("E13").Value = "  -2.55%  "
$ws.Range("E14").Value = "  -2.88%  "
$ws.Range("E15").Value = "  -3.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.71"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.598.76"
$ws.Range("E17").Value = "  -1.85%  "
$ws.Range("E18").Value = "  -2.18%  "
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "208.11"
$ws.Range("E20").Value = "  -4.11%  "
$ws.Range("E21").Value = "  -3.30%  "
$ws.Range("E22").Value = "  -2.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.33"
$ws.Range("E23").Value = "  -3.71%  "
$ws.Range("E24").Value = "  -2.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.78"
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E27").Value = "  -0.83%  "
$ws.Range("E28").Value = "  -3.44%  "
$ws.Range("E29").Value = "  -2.20%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("E31").Value = "  -2.10%  "
$ws.Range("E32").Value = "  -4.36%  "
$ws.Range("E33").Value = "  +20.43%  "
$ws.Range("E34").Value = "  -2.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.306.15"
$ws.Range("E35").Value = "  -3.15%  "
$ws.Range("E36").Value = "  -1.46%  "
$ws.Range("E37").Value = "  -5.61%  "
$ws.Range("E38").Value = "  -3.27%  "
$ws.Range("E39").Value = "  -3.46%  "
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("E41").Value = "  -1.29%  "
$ws.Range("E42").Value = "  +2.48%  "
$ws.Range("E43").Value = "  -3.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.61"
$ws.Range("E44").Value = "  -4.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.723.67"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.50"
$ws.Range("E46").Value = "  -1.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.60"
$ws.Range("E47").Value = "  -1.17%  "
$ws.Range("E48").Value = "  -2.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0505"
$ws.Range("E49").Value = "  -1.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0976"
$ws.Range("E50").Value = "  -1.69%  "
$ws.Range("E51").Value = "  -1.81%  "
